# Update "想去人数" (F column) counts across sheets to match the refreshed
# scrape output (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 219
$ws1.Range("F3").Value  = 54827
$ws1.Range("F9").Value  = 411
$ws1.Range("F10").Value = 3084
$ws1.Range("F14").Value = 1082
$ws1.Range("F20").Value = 108
$ws1.Range("F29").Value = 5160
$ws1.Range("F31").Value = 5088
$ws1.Range("F32").Value = 9094
$ws1.Range("F37").Value = 439
$ws1.Range("F40").Value = 4227

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 1140

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 575

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 411
$ws4.Range("F8").Value  = 3084
$ws4.Range("F14").Value = 1082
$ws4.Range("F21").Value = 108
$ws4.Range("F28").Value = 5160
$ws4.Range("F30").Value = 9094
$ws4.Range("F36").Value = 439
$ws4.Range("F41").Value = 4227
